$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.937.01'
$ws.Range('E2').Value = '  +2.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.652.32'
$ws.Range('E3').Value = '  +2.95%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.87'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.512'
$ws.Range('E6').Value = '  +2.52%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.251'
$ws.Range('E8').Value = '  +2.70%  '
$ws.Range('E9').Value = '  +1.59%  '
$ws.Range('E10').Value = '  +5.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0875'
$ws.Range('E11').Value = '  +2.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.886.90'
$ws.Range('E12').Value = '  +3.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.654.59'
$ws.Range('E13').Value = '  +2.97%  '
$ws.Range('E14').Value = '  +1.92%  '
$ws.Range('E15').Value = '  +2.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.12'
$ws.Range('E16').Value = '  +2.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.946.31'
$ws.Range('E17').Value = '  +2.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '235.87'
$ws.Range('E18').Value = '  +2.45%  '
$ws.Range('E19').Value = '  +1.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.71'
$ws.Range('E20').Value = '  +0.54%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('E22').Value = '  +3.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.30'
$ws.Range('E23').Value = '  +3.85%  '
$ws.Range('E24').Value = '  +3.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.33'
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.10'
$ws.Range('E26').Value = '  +2.08%  '
$ws.Range('E27').Value = '  +0.99%  '
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('E29').Value = '  +2.45%  '
$ws.Range('E30').Value = '  +0.48%  '
$ws.Range('E31').Value = '  +1.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.544.54'
$ws.Range('E32').Value = '  +3.94%  '
$ws.Range('E33').Value = '  +2.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.08'
$ws.Range('E34').Value = '  +5.10%  '
$ws.Range('E35').Value = '  +9.09%  '
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('E37').Value = '  +4.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.894'
$ws.Range('E38').Value = '  +8.94%  '
$ws.Range('E39').Value = '  +2.99%  '
$ws.Range('E40').Value = '  +3.31%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.24'
$ws.Range('E42').Value = '  +2.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.63'
$ws.Range('E43').Value = '  +7.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.793.48'
$ws.Range('E44').Value = '  +2.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.774'
$ws.Range('E45').Value = '  +1.84%  '
$ws.Range('E46').Value = '  -0.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.42'
$ws.Range('E47').Value = '  +1.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.52'
$ws.Range('E48').Value = '  +1.93%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0990'
$ws.Range('E49').Value = '  +3.19%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0505'
$ws.Range('E50').Value = '  +0.96%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.63'
$ws.Range('E51').Value = '  +3.26%  '
